# Objectives.xlsx - saved progress
#
# Row 7 (Ref "DEV-GV-300", "Extend pattern") previously had its Status
# column (D) set to "Dev". Update it to "Val" and record that it has been
# merged on DEV in a new column (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Val"
$ws.Range("E7").Value = "Merged on DEV"

# Leave the selection where the user last clicked while making this edit.
$ws.Range("B5").Select()
